# Hindalco prices sheet update: insert a new latest-date row at the top (row 2),
# shifting all existing data rows down by one. The new row duplicates the
# previous top row's Basic Price / Circular Date / Circular Link (no new
# circular yet) but carries the new Date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the values currently in row 2 (the newest existing row) ---
$origB2 = $ws.Range("B2").Value()
$origC2 = $ws.Range("C2").Value()
$origD2 = $ws.Range("D2").Value()
$origE2 = $ws.Range("E2").Value()
$origF2 = $ws.Range("F2").Value()

# --- insert a new blank row above row 2, shifting everything else down ---
$ws.Rows.Item(2).Insert()

# --- populate the new row 2 with the new date, re-using the other columns ---
$ws.Range("A2").Value = "11-02-2026"
$ws.Range("B2").Value = $origB2
$ws.Range("C2").Value = $origC2
$ws.Range("D2").Value = $origD2
$ws.Range("E2").Value = $origE2
$ws.Range("F2").Value = $origF2

# --- figure out how many data rows now exist ---
$lastRow = $ws.UsedRange.Rows.Count

# --- hyperlinks do not follow the row-insert shift in this engine, so
#     clear every hyperlink on the sheet and rebuild them from the
#     (now-correct) text already sitting in column F of each row ---
$ws.Range("F2").Hyperlinks.Delete()

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $url = $cell.Value()
    if ($url) {
        $ws.Hyperlinks.Add($cell, $url) | Out-Null
    }
}

# --- dimension ref is maintained automatically by the engine on save ---
